$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-15, columns A-F (Giorno, TISG, fcs, buy, MYDIR, need_to_buy)
$data = @(
    @(46036, 12633.9620460958, 11944.3517633521, 20499.86, 8286.27974142422, -11.2178539676546),
    @(46037, 12394.1867748432, 11729.5839694007, 12075.86, 8439.03320374433, 337.198215547709),
    @(46038, 12141.1192610206, 10828.7841653099, 12075.86, 8267.34846907635, 292.511359766096),
    @(46039, 4614.8387139075, 7288.66853077458, 12075.86, 7714.66460363798, 121.97804726719),
    @(46040, 4473.59977441758, 7040.9316671362, 12075.86, 7625.54860000426, 107.942511130852),
    @(46041, 12733.0845017406, 11340.0896915952, 12075.86, 8656.56613530553, 330.033159454196),
    @(46042, 12733.0845017406, 11958.5579263708, 12075.86, 8656.56613530553, 355.802669236512),
    @(46043, 12733.0845017406, 12373.0502944011, 12075.86, 8656.56613530553, 373.07318457111),
    @(46044, 12733.0845017406, 12204.3374440644, 12075.86, 8656.56613530553, 366.043482473745),
    @(46045, 12733.0845017406, 11568.9259833266, 12075.86, 8656.56613530553, 339.568004943007),
    @(46046, 5107.54318705847, 8176.92032683874, 12075.86, 8263.42017596073, 181.853354283311),
    @(46047, 5000.01932310789, 8183.95280080783, 12075.86, 8255.67562535907, 181.823684423621),
    @(46048, 12315.9682835607, 12140.4441206036, 12075.86, 8313.25794670729, 349.076752804619),
    @(46049, 12315.9682835607, 12560.8138681605, 12075.86, 8313.25794670729, 366.592158952827)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}
